$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) for new rows 50 and 51 in column A, copying from row 49
$ws.Cells.Item(49, 1).Copy()
$ws.Cells.Item(50, 1).PasteSpecial(-4122)
$ws.Cells.Item(51, 1).PasteSpecial(-4122)

# Update data values for rows 2 through 51, columns A-E
$ws.Cells.Item(2, 1).Value = 26
$ws.Cells.Item(2, 2).Value = 0.22481
$ws.Cells.Item(2, 3).Value = 0.2245190295220913
$ws.Cells.Item(2, 4).Value = 0.0002909704779086864
$ws.Cells.Item(2, 5).Value = 0.9997091141612829
$ws.Cells.Item(3, 1).Value = 29
$ws.Cells.Item(3, 2).Value = 0.15737
$ws.Cells.Item(3, 3).Value = 0.156958355040731
$ws.Cells.Item(3, 4).Value = 0.0004116449592690141
$ws.Cells.Item(3, 5).Value = 0.9995885244225783
$ws.Cells.Item(4, 1).Value = 13
$ws.Cells.Item(4, 2).Value = 0.17142
$ws.Cells.Item(4, 3).Value = 0.1721887178088827
$ws.Cells.Item(4, 4).Value = 0.0007687178088826851
$ws.Cells.Item(4, 5).Value = 0.9992318726642798
$ws.Cells.Item(5, 1).Value = 45
$ws.Cells.Item(5, 2).Value = 0.18125
$ws.Cells.Item(5, 3).Value = 0.1797307402761048
$ws.Cells.Item(5, 4).Value = 0.001519259723895244
$ws.Cells.Item(5, 5).Value = 0.9984830449248533
$ws.Cells.Item(6, 1).Value = 35
$ws.Cells.Item(6, 2).Value = 0.22368
$ws.Cells.Item(6, 3).Value = 0.2217398189483024
$ws.Cells.Item(6, 4).Value = 0.001940181051697576
$ws.Cells.Item(6, 5).Value = 0.9980635759615301
$ws.Cells.Item(7, 1).Value = 10
$ws.Cells.Item(7, 2).Value = 0.16776
$ws.Cells.Item(7, 3).Value = 0.1652760392406877
$ws.Cells.Item(7, 4).Value = 0.002483960759312315
$ws.Cells.Item(7, 5).Value = 0.9975221940135273
$ws.Cells.Item(8, 1).Value = 18
$ws.Cells.Item(8, 2).Value = 0.16776
$ws.Cells.Item(8, 3).Value = 0.1652760392406877
$ws.Cells.Item(8, 4).Value = 0.002483960759312315
$ws.Cells.Item(8, 5).Value = 0.9975221940135273
$ws.Cells.Item(9, 1).Value = 12
$ws.Cells.Item(9, 2).Value = 0.15877
$ws.Cells.Item(9, 3).Value = 0.156094703907413
$ws.Cells.Item(9, 4).Value = 0.00267529609258696
$ws.Cells.Item(9, 5).Value = 0.9973318420200313
$ws.Cells.Item(10, 1).Value = 24
$ws.Cells.Item(10, 2).Value = 0.18266
$ws.Cells.Item(10, 3).Value = 0.1791100354668042
$ws.Cells.Item(10, 4).Value = 0.003549964533195837
$ws.Cells.Item(10, 5).Value = 0.9964625931357119
$ws.Cells.Item(11, 1).Value = 32
$ws.Cells.Item(11, 2).Value = 0.16186
$ws.Cells.Item(11, 3).Value = 0.1654959149585054
$ws.Cells.Item(11, 4).Value = 0.003635914958505376
$ws.Cells.Item(11, 5).Value = 0.9963772570268614
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = 0.16439
$ws.Cells.Item(12, 3).Value = 0.1682189277482268
$ws.Cells.Item(12, 4).Value = 0.003828927748226796
$ws.Cells.Item(12, 5).Value = 0.9961856770188763
$ws.Cells.Item(13, 1).Value = 28
$ws.Cells.Item(13, 2).Value = 0.16018
$ws.Cells.Item(13, 3).Value = 0.1645314940531793
$ws.Cells.Item(13, 4).Value = 0.004351494053179333
$ws.Cells.Item(13, 5).Value = 0.9956673594065975
$ws.Cells.Item(14, 1).Value = 25
$ws.Cells.Item(14, 2).Value = 0.15315
$ws.Cells.Item(14, 3).Value = 0.1576436960747196
$ws.Cells.Item(14, 4).Value = 0.004493696074719611
$ws.Cells.Item(14, 5).Value = 0.9955264068930648
$ws.Cells.Item(15, 1).Value = 20
$ws.Cells.Item(15, 2).Value = 0.15315
$ws.Cells.Item(15, 3).Value = 0.1576436960747196
$ws.Cells.Item(15, 4).Value = 0.004493696074719611
$ws.Cells.Item(15, 5).Value = 0.9955264068930648
$ws.Cells.Item(16, 1).Value = 44
$ws.Cells.Item(16, 2).Value = 0.15315
$ws.Cells.Item(16, 3).Value = 0.1576436960747196
$ws.Cells.Item(16, 4).Value = 0.004493696074719611
$ws.Cells.Item(16, 5).Value = 0.9955264068930648
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 0.15315
$ws.Cells.Item(17, 3).Value = 0.1576436960747196
$ws.Cells.Item(17, 4).Value = 0.004493696074719611
$ws.Cells.Item(17, 5).Value = 0.9955264068930648
$ws.Cells.Item(18, 1).Value = 33
$ws.Cells.Item(18, 2).Value = 0.16748
$ws.Cells.Item(18, 3).Value = 0.162933511863766
$ws.Cells.Item(18, 4).Value = 0.004546488136233984
$ws.Cells.Item(18, 5).Value = 0.9954740888650466
$ws.Cells.Item(19, 1).Value = 21
$ws.Cells.Item(19, 2).Value = 0.16861
$ws.Cells.Item(19, 3).Value = 0.1732360304203135
$ws.Cells.Item(19, 4).Value = 0.004626030420313537
$ws.Cells.Item(19, 5).Value = 0.9953952711952146
$ws.Cells.Item(20, 1).Value = 8
$ws.Cells.Item(20, 2).Value = 0.16158
$ws.Cells.Item(20, 3).Value = 0.1665477079343512
$ws.Cells.Item(20, 4).Value = 0.004967707934351223
$ws.Cells.Item(20, 5).Value = 0.9950568482000659
$ws.Cells.Item(21, 1).Value = 14
$ws.Cells.Item(21, 2).Value = 0.16158
$ws.Cells.Item(21, 3).Value = 0.1666916511205784
$ws.Cells.Item(21, 4).Value = 0.005111651120578375
$ws.Cells.Item(21, 5).Value = 0.994914344973636
$ws.Cells.Item(22, 1).Value = 36
$ws.Cells.Item(22, 2).Value = 0.20711
$ws.Cells.Item(22, 3).Value = 0.2124249656365614
$ws.Cells.Item(22, 4).Value = 0.005314965636561447
$ws.Cells.Item(22, 5).Value = 0.9947131338752169
$ws.Cells.Item(23, 1).Value = 5
$ws.Cells.Item(23, 2).Value = 0.22762
$ws.Cells.Item(23, 3).Value = 0.2219268355261723
$ws.Cells.Item(23, 4).Value = 0.005693164473827672
$ws.Cells.Item(23, 5).Value = 0.9943390641649571
$ws.Cells.Item(24, 1).Value = 43
$ws.Cells.Item(24, 2).Value = 0.20682
$ws.Cells.Item(24, 3).Value = 0.2129486219422781
$ws.Cells.Item(24, 4).Value = 0.00612862194227809
$ws.Cells.Item(24, 5).Value = 0.9939087092757116
$ws.Cells.Item(25, 1).Value = 41
$ws.Cells.Item(25, 2).Value = 0.18462
$ws.Cells.Item(25, 3).Value = 0.1783048300408952
$ws.Cells.Item(25, 4).Value = 0.006315169959104805
$ws.Cells.Item(25, 5).Value = 0.99372446113541
$ws.Cells.Item(26, 1).Value = 48
$ws.Cells.Item(26, 2).Value = 0.17142
$ws.Cells.Item(26, 3).Value = 0.1648298153633798
$ws.Cells.Item(26, 4).Value = 0.0065901846366202
$ws.Cells.Item(26, 5).Value = 0.9934529615555516
$ws.Cells.Item(27, 1).Value = 19
$ws.Cells.Item(27, 2).Value = 0.17732
$ws.Cells.Item(27, 3).Value = 0.170629388327448
$ws.Cells.Item(27, 4).Value = 0.006690611672552027
$ws.Cells.Item(27, 5).Value = 0.9933538551020795
$ws.Cells.Item(28, 1).Value = 15
$ws.Cells.Item(28, 2).Value = 0.16861
$ws.Cells.Item(28, 3).Value = 0.1753306556431753
$ws.Cells.Item(28, 4).Value = 0.006720655643175283
$ws.Cells.Item(28, 5).Value = 0.9933242100422767
$ws.Cells.Item(29, 1).Value = 37
$ws.Cells.Item(29, 2).Value = 0.17029
$ws.Cells.Item(29, 3).Value = 0.1633754070507057
$ws.Cells.Item(29, 4).Value = 0.006914592949294335
$ws.Cells.Item(29, 5).Value = 0.9931328903188886
$ws.Cells.Item(30, 1).Value = 30
$ws.Cells.Item(30, 2).Value = 0.1939
$ws.Cells.Item(30, 3).Value = 0.2008229979982367
$ws.Cells.Item(30, 4).Value = 0.006922997998236674
$ws.Cells.Item(30, 5).Value = 0.9931246003795727
$ws.Cells.Item(31, 1).Value = 40
$ws.Cells.Item(31, 2).Value = 0.18378
$ws.Cells.Item(31, 3).Value = 0.1768402465357402
$ws.Cells.Item(31, 4).Value = 0.00693975346425979
$ws.Cells.Item(31, 5).Value = 0.9931080747975393
$ws.Cells.Item(32, 1).Value = 39
$ws.Cells.Item(32, 2).Value = 0.15315
$ws.Cells.Item(32, 3).Value = 0.1607856339090119
$ws.Cells.Item(32, 4).Value = 0.007635633909011924
$ws.Cells.Item(32, 5).Value = 0.9924222271900107
$ws.Cells.Item(33, 1).Value = 23
$ws.Cells.Item(33, 2).Value = 0.17732
$ws.Cells.Item(33, 3).Value = 0.1695964700346415
$ws.Cells.Item(33, 4).Value = 0.007723529965358494
$ws.Cells.Item(33, 5).Value = 0.9923356657498869
$ws.Cells.Item(34, 1).Value = 9
$ws.Cells.Item(34, 2).Value = 0.17535
$ws.Cells.Item(34, 3).Value = 0.1669551346146746
$ws.Cells.Item(34, 4).Value = 0.008394865385325356
$ws.Cells.Item(34, 5).Value = 0.9916750216869484
$ws.Cells.Item(35, 1).Value = 42
$ws.Cells.Item(35, 2).Value = 0.17535
$ws.Cells.Item(35, 3).Value = 0.1668399800656924
$ws.Cells.Item(35, 4).Value = 0.008510019934307639
$ws.Cells.Item(35, 5).Value = 0.9915617894060568
$ws.Cells.Item(36, 1).Value = 46
$ws.Cells.Item(36, 2).Value = 0.20345
$ws.Cells.Item(36, 3).Value = 0.2119668834490189
$ws.Cells.Item(36, 4).Value = 0.008516883449018908
$ws.Cells.Item(36, 5).Value = 0.9915550412801299
$ws.Cells.Item(37, 1).Value = 7
$ws.Cells.Item(37, 2).Value = 0.17704
$ws.Cells.Item(37, 3).Value = 0.1677625419157248
$ws.Cells.Item(37, 4).Value = 0.009277458084275242
$ws.Cells.Item(37, 5).Value = 0.990807821962174
$ws.Cells.Item(38, 1).Value = 0
$ws.Cells.Item(38, 2).Value = 0.19109
$ws.Cells.Item(38, 3).Value = 0.181549078802982
$ws.Cells.Item(38, 4).Value = 0.009540921197018049
$ws.Cells.Item(38, 5).Value = 0.9905492476860617
$ws.Cells.Item(39, 1).Value = 27
$ws.Cells.Item(39, 2).Value = 0.17338
$ws.Cells.Item(39, 3).Value = 0.1633322240948365
$ws.Cells.Item(39, 4).Value = 0.01004777590516348
$ws.Cells.Item(39, 5).Value = 0.9900521775852048
$ws.Cells.Item(40, 1).Value = 11
$ws.Cells.Item(40, 2).Value = 0.17985
$ws.Cells.Item(40, 3).Value = 0.1899529192201625
$ws.Cells.Item(40, 4).Value = 0.01010291922016249
$ws.Cells.Item(40, 5).Value = 0.9899981288758553
$ws.Cells.Item(41, 1).Value = 22
$ws.Cells.Item(41, 2).Value = 0.1672
$ws.Cells.Item(41, 3).Value = 0.1781685886344997
$ws.Cells.Item(41, 4).Value = 0.01096858863449973
$ws.Cells.Item(41, 5).Value = 0.9891504159893684
$ws.Cells.Item(42, 1).Value = 6
$ws.Cells.Item(42, 2).Value = 0.21638
$ws.Cells.Item(42, 3).Value = 0.2275050562956946
$ws.Cells.Item(42, 4).Value = 0.01112505629569463
$ws.Cells.Item(42, 5).Value = 0.9889973488181059
$ws.Cells.Item(43, 1).Value = 1
$ws.Cells.Item(43, 2).Value = 0.22368
$ws.Cells.Item(43, 3).Value = 0.2115184983252575
$ws.Cells.Item(43, 4).Value = 0.01216150167474248
$ws.Cells.Item(43, 5).Value = 0.9879846233485271
$ws.Cells.Item(44, 1).Value = 49
$ws.Cells.Item(44, 2).Value = 0.21244
$ws.Cells.Item(44, 3).Value = 0.1997260920445284
$ws.Cells.Item(44, 4).Value = 0.01271390795547159
$ws.Cells.Item(44, 5).Value = 0.9874457061805942
$ws.Cells.Item(45, 1).Value = 47
$ws.Cells.Item(45, 2).Value = 0.18153
$ws.Cells.Item(45, 3).Value = 0.1632332342129318
$ws.Cells.Item(45, 4).Value = 0.01829676578706818
$ws.Cells.Item(45, 5).Value = 0.9820319906712792
$ws.Cells.Item(46, 1).Value = 38
$ws.Cells.Item(46, 2).Value = 0.18153
$ws.Cells.Item(46, 3).Value = 0.1628143091683592
$ws.Cells.Item(46, 4).Value = 0.01871569083164085
$ws.Cells.Item(46, 5).Value = 0.981628151013987
$ws.Cells.Item(47, 1).Value = 31
$ws.Cells.Item(47, 2).Value = 0.22902
$ws.Cells.Item(47, 3).Value = 0.2088530628100272
$ws.Cells.Item(47, 4).Value = 0.02016693718997278
$ws.Cells.Item(47, 5).Value = 0.9802317283036811
$ws.Cells.Item(48, 1).Value = 2
$ws.Cells.Item(48, 2).Value = 0.18856
$ws.Cells.Item(48, 3).Value = 0.1652760392406877
$ws.Cells.Item(48, 4).Value = 0.02328396075931233
$ws.Cells.Item(48, 5).Value = 0.9772458460679528
$ws.Cells.Item(49, 1).Value = 17
$ws.Cells.Item(49, 2).Value = 0.15456
$ws.Cells.Item(49, 3).Value = 0.1238428960053107
$ws.Cells.Item(49, 4).Value = 0.03071710399468927
$ws.Cells.Item(49, 5).Value = 0.970198317389281
$ws.Cells.Item(50, 1).Value = 3
$ws.Cells.Item(50, 2).Value = 0.11606
$ws.Cells.Item(50, 3).Value = 0.1748567354017547
$ws.Cells.Item(50, 4).Value = 0.05879673540175467
$ws.Cells.Item(50, 5).Value = 0.9444683446445983
$ws.Cells.Item(51, 1).Value = 34
$ws.Cells.Item(51, 2).Value = 0.1377
$ws.Cells.Item(51, 3).Value = 0.03279153328981071
$ws.Cells.Item(51, 4).Value = 0.1049084667101893
$ws.Cells.Item(51, 5).Value = 0.9050523460802604
